# Turn the literal powers-of-two in column B into a simple "=A*2" formula
# series (a plotted/derived data column), mirroring the pattern already
# used for the bottom of the column (B15/B16 previously held =B14*2 /
# =B15*2 shared formulas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 is a standalone formula; B3:B21 is written as a single range formula
# so the engine stores it as one contiguous (shared) formula block.
$ws.Range("B2").Formula = "=A2*2"
$ws.Range("B3:B21").Formula = "=A3*2"

# The saved view's active cell/selection moved from F16 to N15.
$ws.Range("N15").Select()
